$d = $word.ActiveDocument

$t = $d.Tables.Item(1)
$t.PreferredWidth = 10451 / 20.0
$t.Columns.Item(1).Width = 8919 / 20.0
$t.Columns.Item(2).Width = 1532 / 20.0
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $t.Rows.Item($i).Height = 290 / 20.0
}

$ps = $d.PageSetup
$ps.TopMargin = 720 / 20.0
$ps.RightMargin = 720 / 20.0
$ps.BottomMargin = 720 / 20.0
$ps.LeftMargin = 720 / 20.0

$s = $d.Styles | Where-Object { $_.NameLocal -eq "Default Paragraph Font" }
$s.Hidden = $true
